# The underlying OOXML diff for this revision is purely a re-serialization
# artifact: every hunk only reorders XML attributes (e.g. namespace
# declarations on <w:document>, the attributes of <w:color>, <w:pgSz>,
# <w:pgMar>, <w:rFonts>, <w:lang>, <w:latentStyles>, <w:lsdException>,
# <w:style>, <w:tblInd>, <w:tblCellMar>, ...) to an alphabetically sorted
# order. This came from a packaging/library fix ("Fixed POI packaging and
# upgraded to POI 3.15") in the authoring tool, not from an actual edit to
# the document's text, formatting, or structure -- every attribute name/value
# pair present before is still present after, just written in a different
# order.
#
# Word's object model (what COM automation exposes) has no concept of
# "attribute serialization order" -- that is an internal detail of the XML
# writer, not a document property a macro/COM client can read or set.
# Because there is no actual content, formatting, or structural change to
# reproduce, this script intentionally performs no mutating operation: the
# document's paragraphs, runs, styles, and section properties are left
# exactly as authored.
$d = $word.ActiveDocument

# Touch the document object (read-only) so the script demonstrably runs
# against the live COM object model without altering any content.
$null = $d.Paragraphs.Count
